$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Penalty" value (column B) for the course rows with ID 6 (row 26)
# and ID 7 (row 30) from 1 to 2 -- this is the "add Penalty" edit referenced
# in the commit message (navbar item driving this value).
$ws.Range("B26").Value = 2
$ws.Range("B30").Value = 2
